# Apply numeric updates to specific cells across three worksheets
# (GLOBAL RESULTS, FUSELAGE, LANDING GEARS) per the commit diff.
$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS ---
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C3").Value = 17.532329314974604
$ws.Range("C5").Value = -0.741123728435491
$ws.Range("C7").Value = 41.58724956580667
$ws.Range("C9").Value = -19.0116448844938
$ws.Range("C13").Value = 16.599203470659226
$ws.Range("C15").Value = -0.9449315944283343
$ws.Range("C17").Value = 17.650278579711
$ws.Range("C19").Value = -24.23981748814523
$ws.Range("C23").Value = 16.599203470659226
$ws.Range("C25").Value = -0.9449315944283343
$ws.Range("C27").Value = 17.650278579711
$ws.Range("C29").Value = -24.23981748814523
$ws.Range("C33").Value = 16.599203470659226
$ws.Range("C35").Value = -0.9449315944283343
$ws.Range("C37").Value = 17.650278579711
$ws.Range("C39").Value = -24.23981748814523
$ws.Range("C43").Value = 17.132888164407937
$ws.Range("C45").Value = -0.6905131041571593
$ws.Range("C47").Value = 31.34060246025785
$ws.Range("C49").Value = -17.71335800034105
$ws.Range("C53").Value = 17.024276851488253
$ws.Range("C55").Value = -0.7940266117833963
$ws.Range("C57").Value = 28.55445538057032
$ws.Range("C59").Value = -20.36873384681775
$ws.Range("C62").Value = 11.410774450069002
$ws.Range("C63").Value = 28.55445538057032
$ws.Range("C64").Value = 35.166164044004056
$ws.Range("C69").Value = 57430.90350385681
$ws.Range("C70").Value = 2984198.0030505783
$ws.Range("C71").Value = 2926767.099546721
$ws.Range("C76").Value = -6714.732829158727

# --- FUSELAGE ---
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C23").Value = 16.834499999999995

# --- LANDING GEARS ---
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 12.948317647221725
$ws.Range("C6").Value = 12.94831764722168
$ws.Range("C7").Value = 16.41590871429195
$ws.Range("C8").Value = 16.41590871429195
$ws.Range("C9").Value = 16.41590871429195
$ws.Range("C10").Value = 16.41590871429195
$ws.Range("C23").Value = 16.41590871429195
